$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 9 (pushes old rows 9..23 down to 11..25),
# carrying the formatting (incl. date style) of row 9 down with them.
$ws.Rows("9:10").Insert()

# New row 9 data
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 45246
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 300000000
$ws.Cells.Item(9, 7).Value = "Espárragos"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 2000
$ws.Cells.Item(9, 12).Value = 2000
$ws.Cells.Item(9, 13).Value = 2000
$ws.Cells.Item(9, 14).Value = "`$/kilo"
$ws.Cells.Item(9, 15).Value = "Provincia de Linares"
$ws.Cells.Item(9, 16).Value = 2000
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# New row 10 data
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 45246
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 300000000
$ws.Cells.Item(10, 7).Value = "Espárragos"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Segunda"
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 1500
$ws.Cells.Item(10, 12).Value = 1500
$ws.Cells.Item(10, 13).Value = 1500
$ws.Cells.Item(10, 14).Value = "`$/kilo"
$ws.Cells.Item(10, 15).Value = "Provincia de Linares"
$ws.Cells.Item(10, 16).Value = 1500
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
